$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.633.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.784.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.15%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.54"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.67%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.87"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.97%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.782.37"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.18%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.30"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.43%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000252"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.95"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.418.95"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.800.18"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.49"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.635.82"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.01"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.19%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.70%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "459.73"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.695"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000152"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +6.07%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.98"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.36%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.03%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.15%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.20"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.62"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.08"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.12%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.34"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.16%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.996"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.74"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.07%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "46.08"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +6.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.24"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.299"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "149.38"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.01%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "392.39"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.97%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.44"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.723.29"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.13%  "
